$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Header row additions (K1:Q1) on "Forecast Comparison" ---
$ws1.Range("K1").Value = "Trend"
$ws1.Range("L1").Value = "Inventory Coverage"
$ws1.Range("M1").Value = "Stockout Risk"
$ws1.Range("N1").Value = "Reorder Urgency"
$ws1.Range("O1").Value = "Sales Trend"
$ws1.Range("P1").Value = "Seasonality Index"
$ws1.Range("Q1").Value = "Lifecycle Stage"

# --- Data rows 2-17, columns K-Q ---
$ws1.Range("K2").Value = "Stable"
$ws1.Range("L2").Value = 6.21
$ws1.Range("M2").Value = "Low"
$ws1.Range("N2").Value = "Normal"
$ws1.Range("O2").Value = "Decreasing (▼)"
$ws1.Range("P2").Value = 0.82
$ws1.Range("Q2").Value = "Decline"

$ws1.Range("K3").Value = "Stable"
$ws1.Range("L3").Value = 5.49
$ws1.Range("M3").Value = "Low"
$ws1.Range("N3").Value = "Normal"
$ws1.Range("O3").Value = "Decreasing (▼)"
$ws1.Range("P3").Value = 0.82
$ws1.Range("Q3").Value = "Decline"

$ws1.Range("K4").Value = "Stable"
$ws1.Range("L4").Value = 4.67
$ws1.Range("M4").Value = "Low"
$ws1.Range("N4").Value = "Normal"
$ws1.Range("O4").Value = "Decreasing (▼)"
$ws1.Range("P4").Value = 0.91
$ws1.Range("Q4").Value = "Decline"

$ws1.Range("K5").Value = "Stable"
$ws1.Range("L5").Value = 3.33
$ws1.Range("M5").Value = "Low"
$ws1.Range("N5").Value = "Normal"
$ws1.Range("O5").Value = "Decreasing (▼)"
$ws1.Range("P5").Value = 0.82
$ws1.Range("Q5").Value = "Decline"

$ws1.Range("K6").Value = "Stable"
$ws1.Range("L6").Value = 2.18
$ws1.Range("M6").Value = "Low"
$ws1.Range("N6").Value = "Normal"
$ws1.Range("O6").Value = "Decreasing (▼)"
$ws1.Range("P6").Value = 1
$ws1.Range("Q6").Value = "Decline"

$ws1.Range("K7").Value = "Stable"
$ws1.Range("L7").Value = 1.19
$ws1.Range("M7").Value = "Low"
$ws1.Range("N7").Value = "Normal"
$ws1.Range("O7").Value = "Decreasing (▼)"
$ws1.Range("P7").Value = 1.1
$ws1.Range("Q7").Value = "Decline"

$ws1.Range("K8").Value = "Stable"
$ws1.Range("L8").Value = 0.2
$ws1.Range("M8").Value = "High"
$ws1.Range("N8").Value = "Urgent"
$ws1.Range("O8").Value = "Decreasing (▼)"
$ws1.Range("P8").Value = 1.07
$ws1.Range("Q8").Value = "Decline"

$ws1.Range("K9").Value = "Stable"
$ws1.Range("L9").Value = 0
$ws1.Range("M9").Value = "High"
$ws1.Range("N9").Value = "Urgent"
$ws1.Range("O9").Value = "Decreasing (▼)"
$ws1.Range("P9").Value = 1.05
$ws1.Range("Q9").Value = "Decline"

$ws1.Range("K10").Value = "Stable"
$ws1.Range("L10").Value = 0
$ws1.Range("M10").Value = "High"
$ws1.Range("N10").Value = "Urgent"
$ws1.Range("O10").Value = "Decreasing (▼)"
$ws1.Range("P10").Value = 0.94
$ws1.Range("Q10").Value = "Decline"

$ws1.Range("K11").Value = "Stable"
$ws1.Range("L11").Value = 0
$ws1.Range("M11").Value = "High"
$ws1.Range("N11").Value = "Urgent"
$ws1.Range("O11").Value = "Decreasing (▼)"
$ws1.Range("P11").Value = 1.06
$ws1.Range("Q11").Value = "Decline"

$ws1.Range("K12").Value = "Stable"
$ws1.Range("L12").Value = 0
$ws1.Range("M12").Value = "High"
$ws1.Range("N12").Value = "Urgent"
$ws1.Range("O12").Value = "Decreasing (▼)"
$ws1.Range("P12").Value = 1
$ws1.Range("Q12").Value = "Decline"

$ws1.Range("K13").Value = "Stable"
$ws1.Range("L13").Value = 0
$ws1.Range("M13").Value = "High"
$ws1.Range("N13").Value = "Urgent"
$ws1.Range("O13").Value = "Decreasing (▼)"
$ws1.Range("P13").Value = 1.1
$ws1.Range("Q13").Value = "Decline"

$ws1.Range("K14").Value = "Stable"
$ws1.Range("L14").Value = 0
$ws1.Range("M14").Value = "High"
$ws1.Range("N14").Value = "Urgent"
$ws1.Range("O14").Value = "Decreasing (▼)"
$ws1.Range("P14").Value = 1.13
$ws1.Range("Q14").Value = "Decline"

$ws1.Range("K15").Value = "Stable"
$ws1.Range("L15").Value = 0
$ws1.Range("M15").Value = "High"
$ws1.Range("N15").Value = "Urgent"
$ws1.Range("O15").Value = "Decreasing (▼)"
$ws1.Range("P15").Value = 0.84
$ws1.Range("Q15").Value = "Decline"

$ws1.Range("K16").Value = "Stable"
$ws1.Range("L16").Value = 0
$ws1.Range("M16").Value = "High"
$ws1.Range("N16").Value = "Urgent"
$ws1.Range("O16").Value = "Decreasing (▼)"
$ws1.Range("P16").Value = 1.1
$ws1.Range("Q16").Value = "Decline"

$ws1.Range("K17").Value = "Stable"
$ws1.Range("L17").Value = 0
$ws1.Range("M17").Value = "High"
$ws1.Range("N17").Value = "Urgent"
$ws1.Range("O17").Value = "Decreasing (▼)"
$ws1.Range("P17").Value = 0.81
$ws1.Range("Q17").Value = "Decline"

# --- Update D7 (MyForecast for week W6) from 112 to 111 ---
$ws1.Range("D7").Value = 111

# --- Update Summary sheet: Max/Min Forecast Week become N/A ---
$ws2.Range("B13").Value = "N/A"
$ws2.Range("B15").Value = "N/A"
